# Perbaiki logika /checkout di routeTransaction.py lalu routeAdmin.py bagian /dataTrans/{id}
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated transaction rows (id_trans, tgl_trans, email_cust, id_staff, status_trans,
# total_harga, metode_byr, tgl_pergi, tgl_balik, id_rute, nama_paket)
$data = @(
    @("A0001", "2024-11-27 03:21:34.900792+07:00", "user2@gmail.com", "ADM_01", "COMPLETED", "1500000.0", "cash",     "2024-12-02", "2024-12-04", "PTKSKW", "Paket Wisata Ketapang"),
    @("A0002", "2024-11-27 03:29:58.357326+07:00", "rio@gmail.com",   "ADM_01", "COMPLETED", "500000.0",  "cash",     "2024-11-30", "2024-12-01", "PTKSKW", "Paket Wisata Singkawang"),
    @("A0003", "2024-11-27 13:48:07.439755+07:00", "user2@gmail.com", "ADM_01", "COMPLETED", "500000.0",  "transfer", "2024-11-30", "2024-12-01", "PTKSKW", "Paket Wisata Singkawang"),
    @("A0004", "2024-11-27 13:50:28.081021+07:00", "rio@gmail.com",   "ADM_01", "COMPLETED", "5220000.0", "cash",     "2024-11-27", "2024-11-29", "PTKSKW", ""),
    @("A0005", "2024-11-27 13:51:53.647278+07:00", "rio@gmail.com",   "KSR_01", "COMPLETED", "8400000.0", "transfer", "2024-11-28", "2024-11-29", "KTPPTK", "")
)

$startRow = 5
$endRow = $startRow + $data.Length - 1

# total_harga (F), tgl_pergi (H) and tgl_balik (I) hold numeric-looking /
# date-looking literal text in the source (e.g. "1500000.0", "2024-12-02").
# Pre-format those columns as Text so the values land as literal strings
# instead of being coerced into numbers / serial dates.
$ws.Range("F$startRow`:F$endRow").NumberFormat = "@"
$ws.Range("H$startRow`:I$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $values[$c]
    }
}
